$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.428.18"
$ws.Range("E2").Value = "  +8.87%  "
$ws.Range("D3").Value = "1.603.59"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +8.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9911"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3683"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3387"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +9.78%  "
$ws.Range("E9").Value = "  +5.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.139"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +7.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07054"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.68"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +8.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.927"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.635"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001086"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +5.21%  "
$ws.Range("D17").Value = "1.601.66"
$ws.Range("E17").Value = "  +8.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9918"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06838"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.03"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +11.23%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.036"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +9.57%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.09"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +10.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.84"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.90%  "
$ws.Range("D24").Value = "22.493.58"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.394"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.537"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +19.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.80"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.59"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +13.06%  "
$ws.Range("D29").Value = "1.782.20"
$ws.Range("E29").Value = "  +8.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.86"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.131"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.058"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +20.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9534"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +15.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08285"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.637"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.265"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +10.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.97"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +14.72%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.623"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +13.38%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.263"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06098"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02225"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +8.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2026"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9915"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5926"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +11.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.833"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.16"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +7.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5720"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +9.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.22"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.986"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +8.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06809"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.87"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +8.86%  "
